$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 3273.5
$ws.Range("I4").Value = 3273.5
$ws.Range("K4").Value = 3273.5
$ws.Range("M4").Value = -3159.5
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("N10").ClearContents()
$ws.Range("H31").Value = 1603.2858
$ws.Range("I31").Value = 1453.8334
$ws.Range("K31").Value = 4361.5002
$ws.Range("M31").Value = -4131.5002
$ws.Range("H62").Value = 3796.4546
$ws.Range("I62").Value = 3686.1
$ws.Range("K62").Value = 3686.1
$ws.Range("M62").Value = -3062.1
$ws.Range("H65").Value = 3796.4546
$ws.Range("I65").Value = 3686.1
$ws.Range("K65").Value = 18430.5
$ws.Range("M65").Value = -15310.5
$ws.Range("I80").Value = 637.7143
$ws.Range("J80").Value = 27778390
$ws.Range("K80").Value = 1913.1429
$ws.Range("L80").Value = 83335170
$ws.Range("M80").Value = -915.1428999999998
$ws.Range("N80").Value = -83337166
$ws.Range("I83").Value = 637.7143
$ws.Range("J83").Value = 27778390
$ws.Range("K83").Value = 5739.428699999999
$ws.Range("L83").Value = 250005510
$ws.Range("M83").Value = -747.4286999999995
$ws.Range("N83").Value = -250015494
$ws.Range("H86").Value = 266669660
$ws.Range("I86").Value = 250003740
$ws.Range("K86").Value = 250003740
$ws.Range("M86").Value = -250002617
$ws.Range("H89").Value = 266669660
$ws.Range("I89").Value = 250003740
$ws.Range("K89").Value = 1250018700
$ws.Range("M89").Value = -1250013084
$ws.Range("H98").Value = 3760.2307
$ws.Range("I98").Value = 2500
$ws.Range("J98").Value = 7961
$ws.Range("K98").Value = 2500
$ws.Range("L98").Value = 7961
$ws.Range("M98").Value = -1002
$ws.Range("N98").Value = -10957
$ws.Range("H122").Value = 3760.2307
$ws.Range("I122").Value = 2500
$ws.Range("J122").Value = 7961
$ws.Range("K122").Value = 7500
$ws.Range("L122").Value = 23883
$ws.Range("M122").Value = -5050
$ws.Range("N122").Value = -28783
$ws.Range("H132").Value = 1585.027
$ws.Range("I132").Value = 1626.3334
$ws.Range("J132").Value = 1408
$ws.Range("K132").Value = 4879.0002
$ws.Range("L132").Value = 4224
$ws.Range("M132").Value = -2349.0002
$ws.Range("N132").Value = -9284
$ws.Range("H137").Value = 2570545.8
$ws.Range("I137").Value = 4546.107
$ws.Range("J137").Value = 9102181
$ws.Range("K137").Value = 13638.321
$ws.Range("L137").Value = 27306543
$ws.Range("M137").Value = -11088.321
$ws.Range("N137").Value = -27311643
$ws.Range("H138").Value = 5206.447
$ws.Range("I138").Value = 8951.65
$ws.Range("K138").Value = 26854.95
$ws.Range("M138").Value = -21714.95

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 11429.714
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 11429.714
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 11429.714
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -12003.714
$ws.Range("H88").Value = 2453
$ws.Range("I88").Value = 2199
$ws.Range("K88").Value = 2199
$ws.Range("M88").Value = -1793
$ws.Range("H91").Value = 2453
$ws.Range("I91").Value = 2199
$ws.Range("K91").Value = 2199
$ws.Range("M91").Value = -795
$ws.Range("H132").Value = 3065.9333
$ws.Range("J132").Value = 4257.5835
$ws.Range("L132").Value = 12772.7505
$ws.Range("N132").Value = -17832.7505

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1175.6923
$ws.Range("I20").Value = 1111.5
$ws.Range("J20").Value = 1389.6666
$ws.Range("K20").Value = 1111.5
$ws.Range("L20").Value = 1389.6666
$ws.Range("M20").Value = -864.5
$ws.Range("N20").Value = -1883.6666
$ws.Range("H95").Value = 42500
$ws.Range("J95").Value = 42500
$ws.Range("L95").Value = 42500
$ws.Range("N95").Value = -47992
$ws.Range("H99").Value = 18881.928
$ws.Range("I99").Value = 20958.6
$ws.Range("K99").Value = 20958.6
$ws.Range("M99").Value = -19460.6
$ws.Range("H105").Value = 21361.273
$ws.Range("J105").Value = 24624.5
$ws.Range("L105").Value = 24624.5
$ws.Range("N105").Value = -28118.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 3655.6365
$ws.Range("J62").Value = 3934.2222
$ws.Range("L62").Value = 3934.2222
$ws.Range("N62").Value = -5182.2222
$ws.Range("H65").Value = 3655.6365
$ws.Range("J65").Value = 3934.2222
$ws.Range("L65").Value = 19671.111
$ws.Range("N65").Value = -25911.111
$ws.Range("H132").Value = 2732
$ws.Range("I132").Value = 3076.2222
$ws.Range("J132").Value = 2422.2
$ws.Range("K132").Value = 9228.6666
$ws.Range("L132").Value = 7266.599999999999
$ws.Range("M132").Value = -6698.6666
$ws.Range("N132").Value = -12326.6
$ws.Range("H135").Value = 114750
$ws.Range("J135").Value = 114750
$ws.Range("L135").Value = 114750
$ws.Range("N135").Value = -124890

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 1131.0278
$ws.Range("I2").Value = 506.85715
$ws.Range("J2").Value = 1528.2273
$ws.Range("K2").Value = 3041.1429
$ws.Range("L2").Value = 9169.363799999999
$ws.Range("M2").Value = -2928.1429
$ws.Range("N2").Value = -9395.363799999999
$ws.Range("H5").Value = 1150.9375
$ws.Range("I5").Value = 881.6
$ws.Range("J5").Value = 1599.8334
$ws.Range("K5").Value = 2644.8
$ws.Range("L5").Value = 4799.5002
$ws.Range("M5").Value = -2532.8
$ws.Range("N5").Value = -5023.5002
$ws.Range("H56").Value = 10005
$ws.Range("I56").Value = 10005
$ws.Range("K56").Value = 10005
$ws.Range("M56").Value = -9475
$ws.Range("H107").Value = 25641828
$ws.Range("J107").Value = 47619864
$ws.Range("L107").Value = 142859592
$ws.Range("N107").Value = -142863432
$ws.Range("H121").Value = 45459756
$ws.Range("J121").Value = 5622.35
$ws.Range("L121").Value = 16867.05
$ws.Range("N121").Value = -19487.05
$ws.Range("H132").Value = 2111.111
$ws.Range("J132").Value = 2333.3333
$ws.Range("L132").Value = 20999.9997
$ws.Range("N132").Value = -26059.9997
$ws.Range("H135").Value = 1150.9375
$ws.Range("I135").Value = 881.6
$ws.Range("J135").Value = 1599.8334
$ws.Range("K135").Value = 7934.400000000001
$ws.Range("L135").Value = 14398.5006
$ws.Range("M135").Value = -5399.400000000001
$ws.Range("N135").Value = -19468.5006
$ws.Range("H137").Value = 2222.077
$ws.Range("J137").Value = 3266.3333
$ws.Range("L137").Value = 9798.999899999999
$ws.Range("N137").Value = -19998.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 40466.75
$ws.Range("J39").Value = 46955.668
$ws.Range("L39").Value = 46955.668
$ws.Range("N39").Value = -48019.668
$ws.Range("H102").Value = 15626190
$ws.Range("I102").Value = 18519428
$ws.Range("K102").Value = 18519428
$ws.Range("M102").Value = -18517806
$ws.Range("H107").Value = 201359.2
$ws.Range("J107").Value = 1999.3334
$ws.Range("L107").Value = 1999.3334
$ws.Range("N107").Value = -5839.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 220205120
$ws.Range("J20").Value = 366666880
$ws.Range("L20").Value = 366666880
$ws.Range("N20").Value = -366667332
$ws.Range("H97").Value = 50000
$ws.Range("J97").Value = 50000
$ws.Range("L97").Value = 50000
$ws.Range("N97").Value = -51982
$ws.Range("H136").Value = 46193.348
$ws.Range("I136").Value = 85928.914
$ws.Range("J136").Value = 2845.4546
$ws.Range("K136").Value = 257786.742
$ws.Range("L136").Value = 8536.363799999999
$ws.Range("M136").Value = -255236.742
$ws.Range("N136").Value = -13636.3638

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 25244
$ws.Range("J96").Value = 68000
$ws.Range("L96").Value = 68000
$ws.Range("N96").Value = -70746
$ws.Range("H122").Value = 1505.4
$ws.Range("I122").Value = 1541.1428
$ws.Range("J122").Value = 1005
$ws.Range("K122").Value = 4623.428400000001
$ws.Range("L122").Value = 3015
$ws.Range("M122").Value = -2173.428400000001
$ws.Range("N122").Value = -7915
$ws.Range("H125").Value = 80176
$ws.Range("J125").Value = 80176
$ws.Range("L125").Value = 80176
$ws.Range("N125").Value = -90016
$ws.Range("H128").Value = 79996.25
$ws.Range("J128").Value = 79996.25
$ws.Range("L128").Value = 79996.25
$ws.Range("N128").Value = -89956.25
$ws.Range("H129").Value = 0
$ws.Range("J129").Value = 0
$ws.Range("L129").Value = 0
$ws.Range("N129").ClearContents()
$ws.Range("H131").Value = 40000
$ws.Range("J131").Value = 40000
$ws.Range("L131").Value = 40000
$ws.Range("N131").Value = -50080
$ws.Range("H132").Value = 2098.087
$ws.Range("I132").Value = 1793.6316
$ws.Range("J132").Value = 3544.25
$ws.Range("K132").Value = 5380.8948
$ws.Range("L132").Value = 10632.75
$ws.Range("M132").Value = -2850.8948
$ws.Range("N132").Value = -15692.75
